$d = $word.ActiveDocument

# The document's "Date" table cell holds "30/10/2019" (the day "30" and the
# rest "/10/2019" are two separate runs, split by a "_GoBack" bookmark).
# Correct the day portion from "30" to "18" -> "18/10/2019", touching only
# that single run of text so the bookmark and surrounding runs are untouched.

$replaced = $false

try {
    $table = $d.Tables.Item(1)
    $dateCell = $table.Cell(3, 2)
    $cellRange = $dateCell.Range
    if ($cellRange.Text -like "*30/10/2019*") {
        $replaced = $cellRange.Find.Execute("30", $true, $false, $false, $false, $false, $true, 1, $false, "18", 2)
    }
} catch {
    $replaced = $false
}

if (-not $replaced) {
    # Fallback: locate the exact whole-word "30" run anywhere in the document.
    $fallback = $d.Content
    $fallback.Find.Execute("30", $true, $false, $false, $false, $false, $true, 1, $false, "18", 2)
}
